$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.394.10"
$ws.Range("E2").Value = "  +3.45%  "

$ws.Range("D3").Value = "1.748.11"
$ws.Range("E3").Value = "  +1.72%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  +0.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.15"
$ws.Range("E5").Value = "  +0.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4816"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2615"
$ws.Range("E8").Value = "  +1.03%  "

$ws.Range("E9").Value = "  -0.26%  "

$ws.Range("D10").Value = "1.746.17"
$ws.Range("E10").Value = "  +1.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.10"
$ws.Range("E11").Value = "  +3.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06937"
$ws.Range("E12").Value = "  -0.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6038"
$ws.Range("E13").Value = "  +0.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.478"
$ws.Range("E14").Value = "  +0.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.35"
$ws.Range("E15").Value = "  +1.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  +0.51%  "

$ws.Range("D17").Value = "27.352.43"
$ws.Range("E17").Value = "  +3.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9989"
$ws.Range("E18").Value = "  +0.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007065"
$ws.Range("E19").Value = "  -0.47%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.49"

$ws.Range("D21").Value = "1.967.81"
$ws.Range("E21").Value = "  +1.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.446"
$ws.Range("E22").Value = "  +1.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.439"
$ws.Range("E23").Value = "  +0.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.111"
$ws.Range("E24").Value = "  +0.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.33"
$ws.Range("E25").Value = "  +3.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.28"
$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.842"
$ws.Range("E27").Value = "  +5.87%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "108.13"
$ws.Range("E28").Value = "  +2.43%  "

$ws.Range("E29").Value = "  -0.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.959"
$ws.Range("E30").Value = "  +1.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07977"
$ws.Range("E31").Value = "  +0.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.673"
$ws.Range("E32").Value = "  +1.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04677"
$ws.Range("E33").Value = "  +4.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.596"
$ws.Range("E34").Value = "  -0.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.014"
$ws.Range("E35").Value = "  +1.73%  "

$ws.Range("E36").Value = "  +0.53%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9249"
$ws.Range("E37").Value = "  -2.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.566"
$ws.Range("E38").Value = "  +7.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.009"
$ws.Range("E39").Value = "  +0.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9992"
$ws.Range("E40").Value = "  +0.38%  "

$ws.Range("E41").Value = "  +5.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01497"
$ws.Range("E42").Value = "  +1.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.86"
$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3843"
$ws.Range("E44").Value = "  +0.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.907"
$ws.Range("E45").Value = "  +0.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1156"
$ws.Range("E46").Value = "  +0.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05360"
$ws.Range("E47").Value = "  +0.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.844"
$ws.Range("E48").Value = "  +1.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.87"
$ws.Range("E49").Value = "  -1.73%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.249"
$ws.Range("E50").Value = "  +3.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.03"
$ws.Range("E51").Value = "  -0.29%  "
